$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (name, count) for rows 2..13
$data = @(
    @("J. Dronkers", 35),
    @("W. C. Ultee", 33),
    @("I. Maas", 33),
    @("P. A. Dykstra", 28),
    @("W. Arts", 27),
    @("J. P. Kleiweg de Zwaan", 26),
    @("P. M. de Graaf", 21),
    @("H. M. Jolles", 20),
    @("G. Kraaykamp", 20),
    @("S. Groenman", 18),
    @("N. Dirk de Graaf", 18),
    @("M. Gesthuizen", 18)
)

# Style reference for column A (bold, bordered) cells, copy from A2 before changing its value
$styleSource = $ws.Range("A2")

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $name = $data[$i][0]
    $count = $data[$i][1]

    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellB = $ws.Cells.Item($rowNum, 2)

    $cellA.Value = $name
    $cellB.Value = $count

    # Ensure the name cell carries the same style as the existing labeled cells
    $styleSource.Copy()
    $cellA.PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0
